$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C3:C7").Value = "N"

$ws.Range("C3").Select()
